$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1 relabeling:
#   E1: "2024e" -> "2024"
#   F1: "2025f" -> "2025e"
# E1's new text ("2024") looks like a number, so force the cell to Text
# format first -- otherwise Excel's normal typing/auto-detection would
# store it as the numeric value 2024 instead of the string "2024".
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2024"

# F1's new text ("2025e") is not a pure number, so it's stored as text as-is.
$ws.Range("F1").Value = "2025e"
